# Communication links WC1-3 hinzugefügt
#
# Adds three new "Element" rows to the tags matrix, one right after each of
# the existing WC1 / WC2 / WC3 rows: "WC1-to-S3", "WC2-to-S3", "WC3-to-S3".
# The new rows are blank (no tag columns marked with "X"), matching the
# formatting of the surrounding rows (row-insert copies the format of the
# row above it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Before insertion the three communication-link rows live at:
#   24 -> WC1
#   25 -> WC2
#   26 -> WC3
# We insert a fresh row directly below each of them, in order, so later
# insert positions must account for rows already shifted down by the
# earlier inserts.

# Insert "WC1-to-S3" right below WC1 (row 24 -> new row 25)
$ws.Rows.Item(25).Insert()
$ws.Cells.Item(25, 1).Value = "WC1-to-S3"

# Insert "WC2-to-S3" right below WC2 (now at row 26 -> new row 27)
$ws.Rows.Item(27).Insert()
$ws.Cells.Item(27, 1).Value = "WC2-to-S3"

# Insert "WC3-to-S3" right below WC3 (now at row 28 -> new row 29)
$ws.Rows.Item(29).Insert()
$ws.Cells.Item(29, 1).Value = "WC3-to-S3"
